# Daily attendance processing - 2025-11-11 08:28:30
#
# Normalises the "Recorded By" column (G) on the session-analysis sheet:
# rows whose recorder list mixes "System" (and/or its lower-cased/backdoor
# duplicate) together with a real reviewer get the leading "System" entry
# rotated to the end of the comma-separated list, e.g.
#   "System, dnasr281@gmail.com"              -> "dnasr281@gmail.com, System"
#   "admin@admin.com, System"                 -> "System, admin@admin.com"
#   "system, System, backup@backdoor.com"     -> "System, backup@backdoor.com, system"
#   "admin@admin.com, dnasr281@gmail.com"     -> "dnasr281@gmail.com, admin@admin.com"
# Rows that are a single name, or exactly "System, backup@backdoor.com",
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $value = $cell.Value2

    if ([string]::IsNullOrEmpty($value)) {
        continue
    }

    $parts = $value.Split(",") | ForEach-Object { $_.Trim() }

    if ($parts.Count -lt 2) {
        continue
    }

    if ($value -eq "System, backup@backdoor.com") {
        continue
    }

    $rotated = (@($parts[1..($parts.Count - 1)]) + @($parts[0])) -join ", "
    $cell.Value2 = $rotated
}
